$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.898.83'
$ws.Range("E2").Value = '  -0.54%  '

$ws.Range("D3").Value = '2.812.23'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.98'
$ws.Range("E5").Value = '  +2.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.26'
$ws.Range("E6").Value = '  -4.20%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.597'
$ws.Range("E9").Value = '  +2.94%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.54'
$ws.Range("E10").Value = '  -5.60%  '

$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.81'
$ws.Range("E13").Value = '  -1.45%  '

$ws.Range("E14").Value = '  +0.21%  '

$ws.Range("D15").Value = '3.253.49'
$ws.Range("E15").Value = '  +0.82%  '

$ws.Range("D16").Value = '2.817.64'
$ws.Range("E16").Value = '  +0.52%  '

$ws.Range("E17").Value = '  +4.25%  '

$ws.Range("D18").Value = '51.736.22'
$ws.Range("E18").Value = '  -0.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.53'
$ws.Range("E19").Value = '  +6.47%  '

$ws.Range("E20").Value = '  -2.99%  '

$ws.Range("E21").Value = '  -0.20%  '

$ws.Range("E22").Value = '  +0.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.26'
$ws.Range("E23").Value = '  +0.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.93'
$ws.Range("E24").Value = '  -0.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.79'
$ws.Range("E25").Value = '  +1.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.82'
$ws.Range("E26").Value = '  +0.75%  '

$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("E28").Value = '  -0.15%  '

$ws.Range("E29").Value = '  +0.48%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0494'
$ws.Range("E30").Value = '  +20.86%  '

$ws.Range("E31").Value = '  -0.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '52.60'
$ws.Range("E32").Value = '  +4.81%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.28'
$ws.Range("E33").Value = '  -1.59%  '

$ws.Range("E34").Value = '  +3.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.48'
$ws.Range("E35").Value = '  +10.29%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0846'
$ws.Range("E36").Value = '  +2.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.09%  '

$ws.Range("E38").Value = '  +0.17%  '

$ws.Range("E39").Value = '  -4.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.26'
$ws.Range("E40").Value = '  -4.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.116'
$ws.Range("E41").Value = '  +0.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '126.77'
$ws.Range("E42").Value = '  -0.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.88'
$ws.Range("E43").Value = '  -2.40%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.48'
$ws.Range("E44").Value = '  -8.11%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.27'
$ws.Range("E45").Value = '  -2.26%  '

$ws.Range("D46").Value = '2.080.81'
$ws.Range("E46").Value = '  +0.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("E47").Value = '  -0.82%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.91'
$ws.Range("E49").Value = '  +6.48%  '

$ws.Range("E50").Value = '  +7.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.07'
$ws.Range("E51").Value = '  +1.08%  '
